# Update WESM exposure data: recomputed ACTUAL_ENERGY, CONTESTABLE_ENERGY,
# TOTAL_BCQ_NOMINATION and WESM_EXPOSURE figures for hours 1-24 (rows 2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 60091.7405
$ws.Range("C2").Value = 5632.2595
$ws.Range("D2").Value = 22500
$ws.Range("E2").Value = 31959.481

$ws.Range("B3").Value = 57606.156
$ws.Range("C3").Value = 5478.844000000001
$ws.Range("D3").Value = 22500
$ws.Range("E3").Value = 29627.31200000001

$ws.Range("B4").Value = 55826.2625
$ws.Range("C4").Value = 5431.737499999999
$ws.Range("D4").Value = 22500
$ws.Range("E4").Value = 27894.52499999999

$ws.Range("B5").Value = 54920.7065
$ws.Range("C5").Value = 5351.2935
$ws.Range("D5").Value = 22500
$ws.Range("E5").Value = 27069.413

$ws.Range("B6").Value = 57738.2435
$ws.Range("C6").Value = 5443.7565
$ws.Range("D6").Value = 22500
$ws.Range("E6").Value = 29794.48699999999

$ws.Range("B7").Value = 62280.3545
$ws.Range("C7").Value = 5521.645500000001
$ws.Range("E7").Value = 34258.709

$ws.Range("B8").Value = 63002.004
$ws.Range("C8").Value = 6568.996000000001
$ws.Range("E8").Value = 33933.008

$ws.Range("B9").Value = 74118.9295
$ws.Range("C9").Value = 7039.0705
$ws.Range("E9").Value = 44579.859

$ws.Range("B10").Value = 89684.016
$ws.Range("C10").Value = 8768.984
$ws.Range("D10").Value = 57000
$ws.Range("E10").Value = 23915.03200000001

$ws.Range("B11").Value = 94141.26149999999
$ws.Range("C11").Value = 14287.7385
$ws.Range("D11").Value = 65000
$ws.Range("E11").Value = 14853.52299999999

$ws.Range("B12").Value = 95535.67999999999
$ws.Range("C12").Value = 16203.32
$ws.Range("D12").Value = 65000
$ws.Range("E12").Value = 14332.35999999999

$ws.Range("B13").Value = 96843.9975
$ws.Range("C13").Value = 16253.0025
$ws.Range("E13").Value = 15590.995

$ws.Range("B14").Value = 96558.0425
$ws.Range("C14").Value = 16193.9575
$ws.Range("E14").Value = 15364.08499999999

$ws.Range("B15").Value = 100667.975
$ws.Range("C15").Value = 16307.025
$ws.Range("E15").Value = 19360.95000000001

$ws.Range("B16").Value = 101388.0945
$ws.Range("C16").Value = 16253.9055
$ws.Range("E16").Value = 20134.18900000001

$ws.Range("B17").Value = 80754.929
$ws.Range("C17").Value = 15667.071
$ws.Range("E17").Value = 87.85800000000745

$ws.Range("B18").Value = 78050.60800000001
$ws.Range("C18").Value = 16597.392
$ws.Range("E18").Value = -3546.783999999992

$ws.Range("C19").Value = 16121.5775
$ws.Range("E19").Value = -3597.647500000006

$ws.Range("C20").Value = 15277.878
$ws.Range("E20").Value = -454.8110000000015

$ws.Range("C21").Value = 13983.214
$ws.Range("E21").Value = -1183.855499999998

$ws.Range("C22").Value = 12008.9865
$ws.Range("E22").Value = 138.1984999999986

$ws.Range("C23").Value = 9648.002000000002
$ws.Range("E23").Value = 1302.546999999991

$ws.Range("C24").Value = 6756.816500000001
$ws.Range("E24").Value = 2128.741999999998

$ws.Range("C25").Value = 6005.93
$ws.Range("E25").Value = 561.2035000000033
